$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column S (year 2022) to the right of the existing data
# (column R, year 2021). Inserting with Shift-right + "format from
# left" copy-origin gives every new S-cell the same cell style as its
# row's R-cell, mirroring the table's existing per-row formatting.
$ws.Range("S1:S33").Insert(-4161, 0) | Out-Null

# Header year.
$ws.Range("S3").Value = 2022

# Men / women / total rates for 2022, row by row (A4:A33 block).
$ws.Range("S4").Value = 5.5
$ws.Range("S5").Value = 8.5
$ws.Range("S6").Value = 2.6

$ws.Range("S7").Value = 16.3
$ws.Range("S8").Value = 25.2
$ws.Range("S9").Value = 7.1

$ws.Range("S10").Value = 1.6
$ws.Range("S11").Value = 3.2
$ws.Range("S12").Value = "-"

$ws.Range("S13").Value = 7.5
$ws.Range("S14").Value = 10.5
$ws.Range("S15").Value = 4.5

$ws.Range("S16").Value = 11.4
$ws.Range("S17").Value = 16.1
$ws.Range("S18").Value = 6.6

$ws.Range("S19").Value = 1.2
$ws.Range("S20").Value = 2.1
$ws.Range("S21").Value = 0.3

$ws.Range("S22").Value = 1.5
$ws.Range("S23").Value = 2.9
$ws.Range("S24").Value = 0

$ws.Range("S25").Value = 0.9
$ws.Range("S26").Value = 1.7
$ws.Range("S27").Value = 0.2

$ws.Range("S28").Value = 14.3
$ws.Range("S29").Value = 22.7
$ws.Range("S30").Value = 7.3

$ws.Range("S31").Value = 1.1
$ws.Range("S32").Value = 2.2
$ws.Range("S33").Value = "-"

$excel.CutCopyMode = $false

# Match the author's final selection state.
$ws.Range("T3").Select() | Out-Null
